# Atualização automática SALDO_PECAS (14/11/2025 20:41)
# Appends one new tracking row to the PRINCIPAL sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162

# Find the first free row right after the current data block (column A).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
$row = $lastRow + 1

$ws.Cells.Item($row, 1).Value  = "DF"
$ws.Cells.Item($row, 2).Value  = "DF00002"
$ws.Cells.Item($row, 3).Value  = ""
$ws.Cells.Item($row, 4).Value  = ""
$ws.Cells.Item($row, 5).Value  = ""
$ws.Cells.Item($row, 6).Value  = "X"
$ws.Cells.Item($row, 7).Value  = "X"
$ws.Cells.Item($row, 8).Value  = "X - (X 03/11/25_12H) - DF"

# Columns I and K hold date-like text ("dd/mm/yy"); force text formatting
# first so Excel doesn't auto-convert the literal strings into date serials,
# then drop the number format again so the cell keeps the sheet's default
# (unstyled) look, matching the rest of the table.
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = "03/11/25"
$ws.Cells.Item($row, 9).ClearFormats()

$ws.Cells.Item($row, 10).Value = "12H"

$ws.Cells.Item($row, 11).NumberFormat = "@"
$ws.Cells.Item($row, 11).Value = "14/11/25"
$ws.Cells.Item($row, 11).ClearFormats()

$ws.Cells.Item($row, 12).Value = "DENTRO"
$ws.Cells.Item($row, 13).Value = ""
